# Auto-generated edit script applying the Valefor_Profits diff
# Sets updated market-price derived columns (H-N) per leve row across
# the ALC, ARM, CUL, LTW and WVR sheets, and clears two cells that the
# diff removes entirely (CUL!M75, CUL!M78).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 418.11365
$ws.Range("J17").Value = 418.11365
$ws.Range("L17").Value = 1254.34095
$ws.Range("N17").Value = -1590.34095
$ws.Range("H76").Value = 3272
$ws.Range("I76").Value = 3390
$ws.Range("J76").Value = 2800
$ws.Range("K76").Value = 3390
$ws.Range("L76").Value = 2800
$ws.Range("M76").Value = -3075
$ws.Range("N76").Value = -3430
$ws.Range("H79").Value = 3272
$ws.Range("I79").Value = 3390
$ws.Range("J79").Value = 2800
$ws.Range("K79").Value = 3390
$ws.Range("L79").Value = 2800
$ws.Range("M79").Value = -2298
$ws.Range("N79").Value = -4984

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 28942.691
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 28942.691
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 28942.691
$ws.Range("N121").Value = -32436.691
$ws.Range("H122").Value = 2194.2683
$ws.Range("I122").Value = 1495.9524
$ws.Range("J122").Value = 2927.5
$ws.Range("K122").Value = 4487.857199999999
$ws.Range("L122").Value = 8782.5
$ws.Range("M122").Value = -2037.857199999999
$ws.Range("N122").Value = -13682.5
$ws.Range("H123").Value = 30000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 30000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
$ws.Range("H124").Value = 14164.5
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 14164.5
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 14164.5
$ws.Range("N124").Value = -23984.5
$ws.Range("H125").Value = 26342.8
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 26342.8
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 26342.8
$ws.Range("N125").Value = -36182.8
$ws.Range("H126").Value = 16666666
$ws.Range("I126").Value = 16666666
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 49999998
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -49997528
$ws.Range("H127").Value = 30000
$ws.Range("I127").Value = 30000
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 30000
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -25040
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H129").Value = 45824
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 45824
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 45824
$ws.Range("N129").Value = -55824
$ws.Range("H130").Value = 33685.8
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 33685.8
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 33685.8
$ws.Range("N130").Value = -43725.8
$ws.Range("H131").Value = 49600
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 49600
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 49600
$ws.Range("N131").Value = -59680
$ws.Range("H132").Value = 1327.4445
$ws.Range("I132").Value = 1085.8334
$ws.Range("J132").Value = 2293.889
$ws.Range("K132").Value = 3257.5002
$ws.Range("L132").Value = 6881.667
$ws.Range("M132").Value = -727.5001999999999
$ws.Range("N132").Value = -11941.667
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H135").Value = 49425
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 49425
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 49425
$ws.Range("N135").Value = -59565
$ws.Range("H137").Value = 30000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 30000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 30000
$ws.Range("N137").Value = -40200
$ws.Range("H138").Value = 41821.75
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 41821.75
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 41821.75
$ws.Range("N138").Value = -52101.75
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 77459.8
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 77459.8
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 77459.8
$ws.Range("N140").Value = -87819.8
$ws.Range("H141").Value = 140429
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 140429
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 140429
$ws.Range("N141").Value = -150789

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 1000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 3000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -4996
$ws.Range("H78").Value = 1000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 9000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -18984
$ws.Range("H101").Value = 8660
$ws.Range("J101").Value = 8660
$ws.Range("L101").Value = 25980
$ws.Range("N101").Value = -30848
$ws.Range("H131").Value = 863.2406999999999
$ws.Range("I131").Value = 377.55554
$ws.Range("J131").Value = 960.3778
$ws.Range("K131").Value = 1132.66662
$ws.Range("L131").Value = 2881.1334
$ws.Range("M131").Value = 3907.33338
$ws.Range("N131").Value = -12961.1334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2024.6818
$ws.Range("I132").Value = 1115.3715
$ws.Range("J132").Value = 5560.8887
$ws.Range("K132").Value = 3346.1145
$ws.Range("L132").Value = 16682.6661
$ws.Range("M132").Value = -816.1144999999997
$ws.Range("N132").Value = -21742.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 19349
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 19349
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 19349
$ws.Range("N119").Value = -29025
$ws.Range("H120").Value = 39000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 39000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 39000
$ws.Range("N120").Value = -48676
$ws.Range("H121").Value = 30971.4
$ws.Range("I121").Value = 30971.4
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 30971.4
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = -29224.4
$ws.Range("H122").Value = 1968.3334
$ws.Range("I122").Value = 1950
$ws.Range("J122").Value = 2005
$ws.Range("K122").Value = 5850
$ws.Range("L122").Value = 6015
$ws.Range("M122").Value = -3400
$ws.Range("N122").Value = -10915
$ws.Range("H123").Value = 21714.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 21714.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 21714.5
$ws.Range("N123").Value = -31514.5
$ws.Range("H124").Value = 27720.666
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 27720.666
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 27720.666
$ws.Range("N124").Value = -37540.666
$ws.Range("H125").Value = 40785.453
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 40785.453
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 40785.453
$ws.Range("N125").Value = -50625.453
$ws.Range("H126").Value = 1237.7059
$ws.Range("I126").Value = 983.5714
$ws.Range("J126").Value = 2423.6667
$ws.Range("K126").Value = 2950.7142
$ws.Range("L126").Value = 7271.000100000001
$ws.Range("M126").Value = -480.7142000000003
$ws.Range("N126").Value = -12211.0001
$ws.Range("H127").Value = 42500
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 42500
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 42500
$ws.Range("N127").Value = -52420
$ws.Range("H128").Value = 35628.75
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 35628.75
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 35628.75
$ws.Range("N128").Value = -45588.75
$ws.Range("H129").Value = 45214.5
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 45214.5
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 45214.5
$ws.Range("N129").Value = -55214.5
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 38566.668
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 38566.668
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 38566.668
$ws.Range("N131").Value = -48646.668
$ws.Range("H132").Value = 1331.3243
$ws.Range("I132").Value = 1309.3334
$ws.Range("J132").Value = 1360.1875
$ws.Range("K132").Value = 3928.0002
$ws.Range("L132").Value = 4080.5625
$ws.Range("M132").Value = -1398.0002
$ws.Range("N132").Value = -9140.5625
$ws.Range("H133").Value = 30715
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 30715
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 30715
$ws.Range("N133").Value = -40835
$ws.Range("H135").Value = 30536.25
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 30536.25
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 30536.25
$ws.Range("N135").Value = -40676.25
$ws.Range("H136").Value = 1153.9302
$ws.Range("I136").Value = 995.6
$ws.Range("J136").Value = 1373.8334
$ws.Range("K136").Value = 2986.8
$ws.Range("L136").Value = 4121.5002
$ws.Range("M136").Value = -436.8000000000002
$ws.Range("N136").Value = -9221.5002
$ws.Range("H137").Value = 33646
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 33646
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 33646
$ws.Range("N137").Value = -43846
$ws.Range("H138").Value = 103194.14
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 103194.14
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 103194.14
$ws.Range("N138").Value = -113474.14
$ws.Range("H139").Value = 37150
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 37150
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 37150
$ws.Range("N139").Value = -47430
$ws.Range("H140").Value = 50764.273
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 50764.273
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 50764.273
$ws.Range("N140").Value = -61124.273
$ws.Range("H141").Value = 33857.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 33857.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 33857.5
$ws.Range("N141").Value = -44217.5

